$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H40" = 3160.8333
    "J40" = 3392.8
    "L40" = 3392.8
    "N40" = -3742.8
    "H64" = 66669710
    "I64" = 200002140
    "K64" = 200002140
    "M64" = -200001892
    "H67" = 66669710
    "I67" = 200002140
    "K67" = 200002140
    "M67" = -200001282
    "H70" = 5114.5713
    "I70" = 5460.4
    "K70" = 16381.2
    "M70" = -16111.2
    "H73" = 5114.5713
    "I73" = 5460.4
    "K73" = 16381.2
    "M73" = -15445.2
    "H74" = 3115.9033
    "I74" = 2418.9092
    "J74" = 3499.25
    "K74" = 2418.9092
    "L74" = 3499.25
    "M74" = -1482.9092
    "N74" = -5371.25
    "H76" = 6811.5454
    "I76" = 5001
    "J76" = 6992.6
    "K76" = 5001
    "L76" = 6992.6
    "M76" = -4686
    "N76" = -7622.6
    "H77" = 3115.9033
    "I77" = 2418.9092
    "J77" = 3499.25
    "K77" = 12094.546
    "L77" = 17496.25
    "M77" = -7414.546
    "N77" = -26856.25
    "H79" = 6811.5454
    "I79" = 5001
    "J79" = 6992.6
    "K79" = 5001
    "L79" = 6992.6
    "M79" = -3909
    "N79" = -9176.6
    "H129" = 857.4286
    "I129" = 500
    "K129" = 1500
    "M129" = 3500
    "H132" = 2545.6667
    "I132" = 1936.5319
    "J132" = 4335
    "K132" = 5809.5957
    "L132" = 13005
    "M132" = -3279.5957
    "N132" = -18065
    "H138" = 1786.26
    "I138" = 941.46155
    "J138" = 3355.1714
    "K138" = 2824.38465
    "L138" = 10065.5142
    "M138" = 2315.61535
    "N138" = -20345.5142
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H5" = 198.57143
    "I5" = 235
    "K5" = 235
    "M5" = -123
    "H10" = 300
    "J10" = 300
    "L10" = 300
    "N10" = -640
    "H32" = 6434.49
    "I32" = 6434.49
    "J32" = 0
    "K32" = 6434.49
    "L32" = 0
    "M32" = -6147.49
    "H63" = 1333.3334
    "I63" = 1333.3334
    "K63" = 1333.3334
    "M63" = -647.3334
    "H66" = 1333.3334
    "I66" = 1333.3334
    "K66" = 6666.666999999999
    "M66" = -3234.666999999999
    "H74" = 2855.4902
    "I74" = 788.86957
    "J74" = 21868.4
    "K74" = 788.86957
    "L74" = 21868.4
    "M74" = 85.13043000000005
    "N74" = -23616.4
    "H77" = 2855.4902
    "I77" = 788.86957
    "J77" = 21868.4
    "K77" = 3944.34785
    "L77" = 109342
    "M77" = 423.6521500000003
    "N77" = -118078
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
$ws.Range("N32").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H4" = 198.57143
    "I4" = 235
    "K4" = 235
    "M4" = -120
    "H22" = 269.25
    "J22" = 384.66666
    "L22" = 384.66666
    "N22" = -730.66666
    "H105" = 1573.9615
    "I105" = 1266.4375
    "K105" = 1266.4375
    "M105" = 480.5625
    "H134" = 1522.9678
    "I134" = 971.8261
    "J134" = 3107.5
    "K134" = 2915.4783
    "L134" = 9322.5
    "M134" = -380.4782999999998
    "N134" = -14392.5
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H5" = 216.23077
    "I5" = 201
    "J5" = 223
    "K5" = 201
    "L5" = 223
    "M5" = -89
    "N5" = -447
    "H31" = 13351340
    "I31" = 45455640
    "J31" = 25027.36
    "K31" = 45455640
    "L31" = 25027.36
    "M31" = -45455345
    "N31" = -25617.36
    "H34" = 13351340
    "I34" = 45455640
    "J34" = 25027.36
    "K34" = 45455640
    "L34" = 25027.36
    "M34" = -45455438
    "N34" = -25431.36
    "H62" = 7890
    "I62" = 9771.429
    "K62" = 9771.429
    "M62" = -9147.429
    "H65" = 7890
    "I65" = 9771.429
    "K65" = 48857.145
    "M65" = -45737.145
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H122" = 931.62164
    "J122" = 1623.0625
    "L122" = 14607.5625
    "N122" = -19507.5625
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H21" = 5000000
    "J21" = 0
    "L21" = 0
    "H30" = 5000000
    "J30" = 0
    "L30" = 0
    "H70" = 3541.6667
    "I70" = 3265.9092
    "K70" = 3265.9092
    "M70" = -2995.9092
    "H73" = 3541.6667
    "I73" = 3265.9092
    "K73" = 3265.9092
    "M73" = -2329.9092
    "H80" = 2471.3
    "I80" = 2002.1666
    "J80" = 3175
    "K80" = 2002.1666
    "L80" = 3175
    "M80" = -1004.1666
    "N80" = -5171
    "H83" = 2471.3
    "I83" = 2002.1666
    "J83" = 3175
    "K83" = 10010.833
    "L83" = 15875
    "M83" = -5018.833000000001
    "N83" = -25859
    "H113" = 1001143.2
    "I113" = 1667605.4
    "J113" = 1450
    "K113" = 1667605.4
    "L113" = 1450
    "M113" = -1665435.4
    "N113" = -5790
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
$ws.Range("N21").ClearContents()
$ws.Range("N30").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H23" = 1002.75
    "I23" = 1002.75
    "K23" = 1002.75
    "M23" = -772.75
    "H68" = 1751.2222
    "I68" = 1521.9524
    "J68" = 2072.2
    "K68" = 1521.9524
    "L68" = 2072.2
    "M68" = -772.9523999999999
    "N68" = -3570.2
    "H71" = 1751.2222
    "I71" = 1521.9524
    "J71" = 2072.2
    "K71" = 7609.762
    "L71" = 10361
    "M71" = -3865.762
    "N71" = -17849
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
